# Updates the crypto price ("Price", column D) and hourly volume change
# ("Volume(1h)", column E) figures on the active sheet to a newer snapshot,
# matching the Sun Sep 22 21:24:05 UTC 2024 GitHub Actions refresh.
#
# Columns D/E are plain text cells (e.g. "63.287.26", "  +0.23%  ") rather
# than numeric cells, so a handful of the new Price values would otherwise
# be auto-coerced into numbers by a plain .Value assignment (e.g. "1.00"
# collapsing to 1, "3.70" collapsing to 3.7). Set-TextValue forces those
# through as literal text (temporarily flips NumberFormat to "@" so Excel
# keeps the string verbatim, then restores the "Normal" style so no stray
# formatting is left behind) while values that are already unambiguous as
# text (contain extra dots, a "%" sign, spaces, ...) are written directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = "63.287.26"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.577.47"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "587.87"
$ws.Range("E5").Value = "  +0.63%  "
Set-TextValue "D6" "144.87"
$ws.Range("E6").Value = "  -1.84%  "
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("E9").Value = "  -2.07%  "
Set-TextValue "D10" "5.61"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("E11").Value = "  -0.12%  "
Set-TextValue "D12" "0.352"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "3.037.95"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "63.192.24"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").Value = "2.576.95"
$ws.Range("E17").Value = "  +1.83%  "
Set-TextValue "D18" "11.08"
$ws.Range("E18").Value = "  -2.37%  "
Set-TextValue "D19" "341.68"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("E20").Value = "  -1.69%  "
Set-TextValue "D21" "6.65"
$ws.Range("E21").Value = "  -3.36%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  +3.74%  "
Set-TextValue "D24" "67.98"
$ws.Range("E24").Value = "  +1.87%  "
Set-TextValue "D25" "1.56"
$ws.Range("E25").Value = "  +5.37%  "
Set-TextValue "D26" "1.62"
$ws.Range("E26").Value = "  -0.59%  "
Set-TextValue "D27" "0.165"
$ws.Range("E27").Value = "  -3.23%  "
Set-TextValue "D28" "1.00"
$ws.Range("E28").Value = "  +0.08%  "
Set-TextValue "D29" "7.95"
$ws.Range("E29").Value = "  -1.70%  "
Set-TextValue "D30" "8.25"
$ws.Range("E30").Value = "  -2.17%  "
Set-TextValue "D31" "1.94"
$ws.Range("E31").Value = "  -2.42%  "
Set-TextValue "D32" "473.39"
$ws.Range("E32").Value = "  +2.53%  "
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("E34").Value = "  +3.93%  "
Set-TextValue "D35" "176.34"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("E36").Value = "  +0.09%  "
Set-TextValue "D37" "0.401"
$ws.Range("E37").Value = "  -1.51%  "
Set-TextValue "D38" "18.90"
$ws.Range("E38").Value = "  -1.58%  "
Set-TextValue "D39" "4.54"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -2.96%  "
Set-TextValue "D42" "40.14"
$ws.Range("E42").Value = "  +1.29%  "
Set-TextValue "D43" "157.73"
$ws.Range("E43").Value = "  +4.37%  "
Set-TextValue "D44" "3.70"
$ws.Range("E44").Value = "  -3.30%  "
Set-TextValue "D45" "21.35"
$ws.Range("E45").Value = "  +2.28%  "
Set-TextValue "D46" "0.635"
$ws.Range("E46").Value = "  +3.57%  "
Set-TextValue "D47" "0.0541"
$ws.Range("E47").Value = "  -1.18%  "
Set-TextValue "D48" "0.0964"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("E49").Value = "  -0.98%  "
Set-TextValue "D50" "18.15"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("E51").Value = "  -0.05%  "
